$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "25.821.41", "0.0780").
# Force each of those cells to text format individually (a comma-joined multi-area
# Range.NumberFormat assignment here only affects the first area) so strings like
# trailing-zero decimals or multi-dot thousand groupings survive as literal text,
# matching the original inline-string cell type instead of being coerced to a number.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D18", "D20", "D21", "D22", "D23", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D38", "D39", "D42", "D43", "D46", "D47")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.821.41"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.635.39"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "215.28"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.0643"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "19.86"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "0.0780"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.25"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.645.14"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "1.860.00"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").Value = "63.13"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "25.820.38"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "4.41"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "194.11"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "9.94"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").Value = "6.15"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "140.08"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "0.0496"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").Value = "3.33"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  +2.43%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "0.553"
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").Value = "1.113.12"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "5.58"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "99.62"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("B46").Value = "SynthetixNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").Value = "  +13.84%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "55.53"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  -5.44%  "
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  +0.35%  "
